$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Jhye Richardson"

# Insert a new column before column A (shifts existing columns A..L to B..M)
$ws.Range("A1").EntireColumn.Insert()

# Insert a new row before row 2 (shifts existing row 2 to row 3)
$ws.Range("A2").EntireRow.Insert()

# The numeric-looking columns (runs, balls, fours, sixes, sr) are stored as
# text in this sheet (see the "numberStoredAsText" ignoredError below), so
# force those cells to Text format before writing, otherwise Excel would
# silently convert them to real numbers.
$ws.Range("E2:I3").NumberFormat = "@"

# Header row (row 1)
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2 - new match (8th)
$ws.Range("A2").Value = "8th"
$ws.Range("B2").Value = "Punjab Kings"
$ws.Range("C2").Value = "Jhye Richardson"
$ws.Range("D2").Value = "b Ali"
$ws.Range("E2").Value = "15"
$ws.Range("F2").Value = "22"
$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "68.18"
$ws.Range("J2").Value = "Chennai Super Kings"
$ws.Range("K2").Value = "Wankhede"
$ws.Range("L2").Value = "April 16"
$ws.Range("M2").Value = "Super Kings won by 6 wickets (with 26 balls remaining)"

# Row 3 - existing match (4th), shifted from row 2 to row 3 and column A to B..M
$ws.Range("A3").Value = "4th"
$ws.Range("B3").Value = "Punjab Kings"
$ws.Range("C3").Value = "Jhye Richardson"
$ws.Range("D3").Value = "c Morris b Sakariya"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "2"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "0.00"
$ws.Range("J3").Value = "Rajasthan Royals"
$ws.Range("K3").Value = "Wankhede"
$ws.Range("L3").Value = "April 12"
$ws.Range("M3").Value = "Punjab Kings won by 4 runs"
